$wb = $excel.ActiveWorkbook
$wsResults = $wb.Worksheets.Item("Results")

# Populate the Monte Carlo / sensitivity analysis results (rows 2-18, columns A-V)
$wsResults.Cells.Item(2,1).Value = 4.9000000000000004
$wsResults.Cells.Item(2,2).Value = 10
$wsResults.Cells.Item(2,3).Value = 11
$wsResults.Cells.Item(2,4).Value = 15
$wsResults.Cells.Item(2,5).Value = 16
$wsResults.Cells.Item(2,6).Value = 1
$wsResults.Cells.Item(2,7).Value = 1
$wsResults.Cells.Item(2,8).Value = 1
$wsResults.Cells.Item(2,9).Value = 1
$wsResults.Cells.Item(2,10).Value = 0.04
$wsResults.Cells.Item(2,11).Value = 0.04
$wsResults.Cells.Item(2,12).Value = 0.04
$wsResults.Cells.Item(2,13).Value = 0.04
$wsResults.Cells.Item(2,14).Value = 50
$wsResults.Cells.Item(2,15).Value = 200
$wsResults.Cells.Item(2,16).Value = 550
$wsResults.Cells.Item(2,17).Value = 0.026778146297450366
$wsResults.Cells.Item(2,18).Value = 0.000043962823000668649
$wsResults.Cells.Item(2,19).Value = 0.75059218247115644
$wsResults.Cells.Item(2,20).Value = 0.00076075678035148722
$wsResults.Cells.Item(2,21).Value = 10000
$wsResults.Cells.Item(2,22).Value = 5

$wsResults.Cells.Item(3,1).Value = 4.9000000000000004
$wsResults.Cells.Item(3,2).Value = 10
$wsResults.Cells.Item(3,3).Value = 11
$wsResults.Cells.Item(3,4).Value = 15
$wsResults.Cells.Item(3,5).Value = 16
$wsResults.Cells.Item(3,6).Value = 1
$wsResults.Cells.Item(3,7).Value = 1
$wsResults.Cells.Item(3,8).Value = 1
$wsResults.Cells.Item(3,9).Value = 1
$wsResults.Cells.Item(3,10).Value = 0.04
$wsResults.Cells.Item(3,11).Value = 0.04
$wsResults.Cells.Item(3,12).Value = 0.04
$wsResults.Cells.Item(3,13).Value = 0.04
$wsResults.Cells.Item(3,14).Value = 50
$wsResults.Cells.Item(3,15).Value = 200
$wsResults.Cells.Item(3,16).Value = 650
$wsResults.Cells.Item(3,17).Value = 0.026359855079758819
$wsResults.Cells.Item(3,18).Value = 0.000065516638209457519
$wsResults.Cells.Item(3,19).Value = 0.75270689500497256
$wsResults.Cells.Item(3,20).Value = 0.001231052177134672
$wsResults.Cells.Item(3,21).Value = 10000
$wsResults.Cells.Item(3,22).Value = 5

$wsResults.Cells.Item(4,1).Value = 4.9000000000000004
$wsResults.Cells.Item(4,2).Value = 10
$wsResults.Cells.Item(4,3).Value = 11
$wsResults.Cells.Item(4,4).Value = 15
$wsResults.Cells.Item(4,5).Value = 16
$wsResults.Cells.Item(4,6).Value = 1
$wsResults.Cells.Item(4,7).Value = 1
$wsResults.Cells.Item(4,8).Value = 1
$wsResults.Cells.Item(4,9).Value = 1
$wsResults.Cells.Item(4,10).Value = 0.04
$wsResults.Cells.Item(4,11).Value = 0.04
$wsResults.Cells.Item(4,12).Value = 0.04
$wsResults.Cells.Item(4,13).Value = 0.04
$wsResults.Cells.Item(4,14).Value = 50
$wsResults.Cells.Item(4,15).Value = 400
$wsResults.Cells.Item(4,16).Value = 550
$wsResults.Cells.Item(4,17).Value = 0.020550519571121802
$wsResults.Cells.Item(4,18).Value = 0.000016755811335919091
$wsResults.Cells.Item(4,19).Value = 0.87568495429540083
$wsResults.Cells.Item(4,20).Value = 0.00078245972291806948
$wsResults.Cells.Item(4,21).Value = 10000
$wsResults.Cells.Item(4,22).Value = 5

$wsResults.Cells.Item(5,1).Value = 4.9000000000000004
$wsResults.Cells.Item(5,2).Value = 10
$wsResults.Cells.Item(5,3).Value = 11
$wsResults.Cells.Item(5,4).Value = 15
$wsResults.Cells.Item(5,5).Value = 16
$wsResults.Cells.Item(5,6).Value = 1
$wsResults.Cells.Item(5,7).Value = 1
$wsResults.Cells.Item(5,8).Value = 1
$wsResults.Cells.Item(5,9).Value = 1
$wsResults.Cells.Item(5,10).Value = 0.04
$wsResults.Cells.Item(5,11).Value = 0.04
$wsResults.Cells.Item(5,12).Value = 0.04
$wsResults.Cells.Item(5,13).Value = 0.04
$wsResults.Cells.Item(5,14).Value = 50
$wsResults.Cells.Item(5,15).Value = 400
$wsResults.Cells.Item(5,16).Value = 650
$wsResults.Cells.Item(5,17).Value = 0.020245401844702294
$wsResults.Cells.Item(5,18).Value = 0.000017154700252746717
$wsResults.Cells.Item(5,19).Value = 0.87497537994888508
$wsResults.Cells.Item(5,20).Value = 0.00045083352932814543
$wsResults.Cells.Item(5,21).Value = 10000
$wsResults.Cells.Item(5,22).Value = 5

$wsResults.Cells.Item(6,1).Value = 4.9000000000000004
$wsResults.Cells.Item(6,2).Value = 10
$wsResults.Cells.Item(6,3).Value = 11
$wsResults.Cells.Item(6,4).Value = 15
$wsResults.Cells.Item(6,5).Value = 16
$wsResults.Cells.Item(6,6).Value = 1
$wsResults.Cells.Item(6,7).Value = 1
$wsResults.Cells.Item(6,8).Value = 1
$wsResults.Cells.Item(6,9).Value = 1
$wsResults.Cells.Item(6,10).Value = 0.04
$wsResults.Cells.Item(6,11).Value = 0.04
$wsResults.Cells.Item(6,12).Value = 0.04
$wsResults.Cells.Item(6,13).Value = 0.04
$wsResults.Cells.Item(6,14).Value = 150
$wsResults.Cells.Item(6,15).Value = 200
$wsResults.Cells.Item(6,16).Value = 550
$wsResults.Cells.Item(6,17).Value = 0.043957818870319264
$wsResults.Cells.Item(6,18).Value = 0.00079814966304870704
$wsResults.Cells.Item(6,19).Value = 0.25194706997238797
$wsResults.Cells.Item(6,20).Value = 0.0061224566195884674
$wsResults.Cells.Item(6,21).Value = 10000
$wsResults.Cells.Item(6,22).Value = 5

$wsResults.Cells.Item(7,1).Value = 4.9000000000000004
$wsResults.Cells.Item(7,2).Value = 10
$wsResults.Cells.Item(7,3).Value = 11
$wsResults.Cells.Item(7,4).Value = 15
$wsResults.Cells.Item(7,5).Value = 16
$wsResults.Cells.Item(7,6).Value = 1
$wsResults.Cells.Item(7,7).Value = 1
$wsResults.Cells.Item(7,8).Value = 1
$wsResults.Cells.Item(7,9).Value = 1
$wsResults.Cells.Item(7,10).Value = 0.04
$wsResults.Cells.Item(7,11).Value = 0.04
$wsResults.Cells.Item(7,12).Value = 0.04
$wsResults.Cells.Item(7,13).Value = 0.04
$wsResults.Cells.Item(7,14).Value = 150
$wsResults.Cells.Item(7,15).Value = 200
$wsResults.Cells.Item(7,16).Value = 650
$wsResults.Cells.Item(7,17).Value = 0.044538518109261356
$wsResults.Cells.Item(7,18).Value = 0.00070767034001242661
$wsResults.Cells.Item(7,19).Value = 0.24598611452073321
$wsResults.Cells.Item(7,20).Value = 0.0065058650103057203
$wsResults.Cells.Item(7,21).Value = 10000
$wsResults.Cells.Item(7,22).Value = 5

$wsResults.Cells.Item(8,1).Value = 4.9000000000000004
$wsResults.Cells.Item(8,2).Value = 10
$wsResults.Cells.Item(8,3).Value = 11
$wsResults.Cells.Item(8,4).Value = 15
$wsResults.Cells.Item(8,5).Value = 16
$wsResults.Cells.Item(8,6).Value = 1
$wsResults.Cells.Item(8,7).Value = 1
$wsResults.Cells.Item(8,8).Value = 1
$wsResults.Cells.Item(8,9).Value = 1
$wsResults.Cells.Item(8,10).Value = 0.04
$wsResults.Cells.Item(8,11).Value = 0.04
$wsResults.Cells.Item(8,12).Value = 0.04
$wsResults.Cells.Item(8,13).Value = 0.04
$wsResults.Cells.Item(8,14).Value = 150
$wsResults.Cells.Item(8,15).Value = 400
$wsResults.Cells.Item(8,16).Value = 550
$wsResults.Cells.Item(8,17).Value = 0.027815436958749711
$wsResults.Cells.Item(8,18).Value = 0.00021681259332896008
$wsResults.Cells.Item(8,19).Value = 0.62929472625482163
$wsResults.Cells.Item(8,20).Value = 0.0018967545717024918
$wsResults.Cells.Item(8,21).Value = 10000
$wsResults.Cells.Item(8,22).Value = 5

$wsResults.Cells.Item(9,1).Value = 4.9000000000000004
$wsResults.Cells.Item(9,2).Value = 10
$wsResults.Cells.Item(9,3).Value = 11
$wsResults.Cells.Item(9,4).Value = 15
$wsResults.Cells.Item(9,5).Value = 16
$wsResults.Cells.Item(9,6).Value = 1
$wsResults.Cells.Item(9,7).Value = 1
$wsResults.Cells.Item(9,8).Value = 1
$wsResults.Cells.Item(9,9).Value = 1
$wsResults.Cells.Item(9,10).Value = 0.04
$wsResults.Cells.Item(9,11).Value = 0.04
$wsResults.Cells.Item(9,12).Value = 0.04
$wsResults.Cells.Item(9,13).Value = 0.04
$wsResults.Cells.Item(9,14).Value = 150
$wsResults.Cells.Item(9,15).Value = 400
$wsResults.Cells.Item(9,16).Value = 650
$wsResults.Cells.Item(9,17).Value = 0.027674836447346463
$wsResults.Cells.Item(9,18).Value = 0.00019194458530667486
$wsResults.Cells.Item(9,19).Value = 0.6245760128162916
$wsResults.Cells.Item(9,20).Value = 0.0015083122669846891
$wsResults.Cells.Item(9,21).Value = 10000
$wsResults.Cells.Item(9,22).Value = 5

$wsResults.Cells.Item(10,1).Value = 4.9000000000000004
$wsResults.Cells.Item(10,2).Value = 10
$wsResults.Cells.Item(10,3).Value = 11
$wsResults.Cells.Item(10,4).Value = 15
$wsResults.Cells.Item(10,5).Value = 16
$wsResults.Cells.Item(10,6).Value = 1
$wsResults.Cells.Item(10,7).Value = 1
$wsResults.Cells.Item(10,8).Value = 1
$wsResults.Cells.Item(10,9).Value = 1
$wsResults.Cells.Item(10,10).Value = 0.02
$wsResults.Cells.Item(10,11).Value = 0.02
$wsResults.Cells.Item(10,12).Value = 0.02
$wsResults.Cells.Item(10,13).Value = 0.02
$wsResults.Cells.Item(10,14).Value = 50
$wsResults.Cells.Item(10,15).Value = 200
$wsResults.Cells.Item(10,16).Value = 550
$wsResults.Cells.Item(10,17).Value = 0.024539692177653512
$wsResults.Cells.Item(10,18).Value = 0.000021646434886623961
$wsResults.Cells.Item(10,19).Value = 0.75037632108335739
$wsResults.Cells.Item(10,20).Value = 0.0013343137874497455
$wsResults.Cells.Item(10,21).Value = 10000
$wsResults.Cells.Item(10,22).Value = 5

$wsResults.Cells.Item(11,1).Value = 4.9000000000000004
$wsResults.Cells.Item(11,2).Value = 10
$wsResults.Cells.Item(11,3).Value = 11
$wsResults.Cells.Item(11,4).Value = 15
$wsResults.Cells.Item(11,5).Value = 16
$wsResults.Cells.Item(11,6).Value = 1
$wsResults.Cells.Item(11,7).Value = 1
$wsResults.Cells.Item(11,8).Value = 1
$wsResults.Cells.Item(11,9).Value = 1
$wsResults.Cells.Item(11,10).Value = 0.02
$wsResults.Cells.Item(11,11).Value = 0.02
$wsResults.Cells.Item(11,12).Value = 0.02
$wsResults.Cells.Item(11,13).Value = 0.02
$wsResults.Cells.Item(11,14).Value = 50
$wsResults.Cells.Item(11,15).Value = 200
$wsResults.Cells.Item(11,16).Value = 650
$wsResults.Cells.Item(11,17).Value = 0.024213497869152573
$wsResults.Cells.Item(11,18).Value = 0.000058184315360275055
$wsResults.Cells.Item(11,19).Value = 0.74992873058871601
$wsResults.Cells.Item(11,20).Value = 0.0016374346347852624
$wsResults.Cells.Item(11,21).Value = 10000
$wsResults.Cells.Item(11,22).Value = 5

$wsResults.Cells.Item(12,1).Value = 4.9000000000000004
$wsResults.Cells.Item(12,2).Value = 10
$wsResults.Cells.Item(12,3).Value = 11
$wsResults.Cells.Item(12,4).Value = 15
$wsResults.Cells.Item(12,5).Value = 16
$wsResults.Cells.Item(12,6).Value = 1
$wsResults.Cells.Item(12,7).Value = 1
$wsResults.Cells.Item(12,8).Value = 1
$wsResults.Cells.Item(12,9).Value = 1
$wsResults.Cells.Item(12,10).Value = 0.02
$wsResults.Cells.Item(12,11).Value = 0.02
$wsResults.Cells.Item(12,12).Value = 0.02
$wsResults.Cells.Item(12,13).Value = 0.02
$wsResults.Cells.Item(12,14).Value = 50
$wsResults.Cells.Item(12,15).Value = 400
$wsResults.Cells.Item(12,16).Value = 550
$wsResults.Cells.Item(12,17).Value = 0.018240823472723662
$wsResults.Cells.Item(12,18).Value = 0.000049349792156898743
$wsResults.Cells.Item(12,19).Value = 0.87563652501947487
$wsResults.Cells.Item(12,20).Value = 0.00085836951990923946
$wsResults.Cells.Item(12,21).Value = 10000
$wsResults.Cells.Item(12,22).Value = 5

$wsResults.Cells.Item(13,1).Value = 4.9000000000000004
$wsResults.Cells.Item(13,2).Value = 10
$wsResults.Cells.Item(13,3).Value = 11
$wsResults.Cells.Item(13,4).Value = 15
$wsResults.Cells.Item(13,5).Value = 16
$wsResults.Cells.Item(13,6).Value = 1
$wsResults.Cells.Item(13,7).Value = 1
$wsResults.Cells.Item(13,8).Value = 1
$wsResults.Cells.Item(13,9).Value = 1
$wsResults.Cells.Item(13,10).Value = 0.02
$wsResults.Cells.Item(13,11).Value = 0.02
$wsResults.Cells.Item(13,12).Value = 0.02
$wsResults.Cells.Item(13,13).Value = 0.02
$wsResults.Cells.Item(13,14).Value = 50
$wsResults.Cells.Item(13,15).Value = 400
$wsResults.Cells.Item(13,16).Value = 650
$wsResults.Cells.Item(13,17).Value = 0.017879609926206443
$wsResults.Cells.Item(13,18).Value = 0.00001183869852412298
$wsResults.Cells.Item(13,19).Value = 0.8763785555182082
$wsResults.Cells.Item(13,20).Value = 0.00049796300225006442
$wsResults.Cells.Item(13,21).Value = 10000
$wsResults.Cells.Item(13,22).Value = 5

$wsResults.Cells.Item(14,1).Value = 4.9000000000000004
$wsResults.Cells.Item(14,2).Value = 10
$wsResults.Cells.Item(14,3).Value = 11
$wsResults.Cells.Item(14,4).Value = 15
$wsResults.Cells.Item(14,5).Value = 16
$wsResults.Cells.Item(14,6).Value = 1
$wsResults.Cells.Item(14,7).Value = 1
$wsResults.Cells.Item(14,8).Value = 1
$wsResults.Cells.Item(14,9).Value = 1
$wsResults.Cells.Item(14,10).Value = 0.02
$wsResults.Cells.Item(14,11).Value = 0.02
$wsResults.Cells.Item(14,12).Value = 0.02
$wsResults.Cells.Item(14,13).Value = 0.02
$wsResults.Cells.Item(14,14).Value = 150
$wsResults.Cells.Item(14,15).Value = 200
$wsResults.Cells.Item(14,16).Value = 550
$wsResults.Cells.Item(14,17).Value = 0.039892509643641903
$wsResults.Cells.Item(14,18).Value = 0.0011825732614869787
$wsResults.Cells.Item(14,19).Value = 0.25628908681405488
$wsResults.Cells.Item(14,20).Value = 0.0069030646042207953
$wsResults.Cells.Item(14,21).Value = 10000
$wsResults.Cells.Item(14,22).Value = 5

$wsResults.Cells.Item(15,1).Value = 4.9000000000000004
$wsResults.Cells.Item(15,2).Value = 10
$wsResults.Cells.Item(15,3).Value = 11
$wsResults.Cells.Item(15,4).Value = 15
$wsResults.Cells.Item(15,5).Value = 16
$wsResults.Cells.Item(15,6).Value = 1
$wsResults.Cells.Item(15,7).Value = 1
$wsResults.Cells.Item(15,8).Value = 1
$wsResults.Cells.Item(15,9).Value = 1
$wsResults.Cells.Item(15,10).Value = 0.02
$wsResults.Cells.Item(15,11).Value = 0.02
$wsResults.Cells.Item(15,12).Value = 0.02
$wsResults.Cells.Item(15,13).Value = 0.02
$wsResults.Cells.Item(15,14).Value = 150
$wsResults.Cells.Item(15,15).Value = 200
$wsResults.Cells.Item(15,16).Value = 650
$wsResults.Cells.Item(15,17).Value = 0.039516862894926383
$wsResults.Cells.Item(15,18).Value = 0.00060546006907180778
$wsResults.Cells.Item(15,19).Value = 0.24213908670634252
$wsResults.Cells.Item(15,20).Value = 0.003793029428375587
$wsResults.Cells.Item(15,21).Value = 10000
$wsResults.Cells.Item(15,22).Value = 5

$wsResults.Cells.Item(16,1).Value = 4.9000000000000004
$wsResults.Cells.Item(16,2).Value = 10
$wsResults.Cells.Item(16,3).Value = 11
$wsResults.Cells.Item(16,4).Value = 15
$wsResults.Cells.Item(16,5).Value = 16
$wsResults.Cells.Item(16,6).Value = 1
$wsResults.Cells.Item(16,7).Value = 1
$wsResults.Cells.Item(16,8).Value = 1
$wsResults.Cells.Item(16,9).Value = 1
$wsResults.Cells.Item(16,10).Value = 0.02
$wsResults.Cells.Item(16,11).Value = 0.02
$wsResults.Cells.Item(16,12).Value = 0.02
$wsResults.Cells.Item(16,13).Value = 0.02
$wsResults.Cells.Item(16,14).Value = 150
$wsResults.Cells.Item(16,15).Value = 400
$wsResults.Cells.Item(16,16).Value = 550
$wsResults.Cells.Item(16,17).Value = 0.022133205741508073
$wsResults.Cells.Item(16,18).Value = 0.000050734586515548537
$wsResults.Cells.Item(16,19).Value = 0.62495074185562127
$wsResults.Cells.Item(16,20).Value = 0.0011058087229678692
$wsResults.Cells.Item(16,21).Value = 10000
$wsResults.Cells.Item(16,22).Value = 5

$wsResults.Cells.Item(17,1).Value = 4.9000000000000004
$wsResults.Cells.Item(17,2).Value = 10
$wsResults.Cells.Item(17,3).Value = 11
$wsResults.Cells.Item(17,4).Value = 15
$wsResults.Cells.Item(17,5).Value = 16
$wsResults.Cells.Item(17,6).Value = 1
$wsResults.Cells.Item(17,7).Value = 1
$wsResults.Cells.Item(17,8).Value = 1
$wsResults.Cells.Item(17,9).Value = 1
$wsResults.Cells.Item(17,10).Value = 0.02
$wsResults.Cells.Item(17,11).Value = 0.02
$wsResults.Cells.Item(17,12).Value = 0.02
$wsResults.Cells.Item(17,13).Value = 0.02
$wsResults.Cells.Item(17,14).Value = 150
$wsResults.Cells.Item(17,15).Value = 400
$wsResults.Cells.Item(17,16).Value = 650
$wsResults.Cells.Item(17,17).Value = 0.021624865597261064
$wsResults.Cells.Item(17,18).Value = 0.00013106489658985156
$wsResults.Cells.Item(17,19).Value = 0.62564164592989302
$wsResults.Cells.Item(17,20).Value = 0.0024818793314518692
$wsResults.Cells.Item(17,21).Value = 10000
$wsResults.Cells.Item(17,22).Value = 5

$wsResults.Cells.Item(18,1).Value = 4.9000000000000004
$wsResults.Cells.Item(18,2).Value = 10
$wsResults.Cells.Item(18,3).Value = 11
$wsResults.Cells.Item(18,4).Value = 15
$wsResults.Cells.Item(18,5).Value = 16
$wsResults.Cells.Item(18,6).Value = 1
$wsResults.Cells.Item(18,7).Value = 1
$wsResults.Cells.Item(18,8).Value = 1
$wsResults.Cells.Item(18,9).Value = 1
$wsResults.Cells.Item(18,10).Value = 0.03
$wsResults.Cells.Item(18,11).Value = 0.03
$wsResults.Cells.Item(18,12).Value = 0.03
$wsResults.Cells.Item(18,13).Value = 0.03
$wsResults.Cells.Item(18,14).Value = 100
$wsResults.Cells.Item(18,15).Value = 300
$wsResults.Cells.Item(18,16).Value = 600
$wsResults.Cells.Item(18,17).Value = 0.023352437180671114
$wsResults.Cells.Item(18,18).Value = 0.000069821935068287361
$wsResults.Cells.Item(18,19).Value = 0.66408766945349051
$wsResults.Cells.Item(18,20).Value = 0.0014568814848845209
$wsResults.Cells.Item(18,21).Value = 10000
$wsResults.Cells.Item(18,22).Value = 5

# Make "Results" the active sheet/tab (was "Test Cases")
$wsResults.Activate()

# Update the selection on the Results sheet to the newly populated data block
$wsResults.Range("A2:V18").Select()

